$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6054.5454
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 6825
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 6825
$ws.Range("M51").Value = -3516
$ws.Range("N51").Value = -7793

$ws.Range("H64").Value = 4180.8335
$ws.Range("I64").Value = 3800
$ws.Range("J64").Value = 4656.875
$ws.Range("K64").Value = 3800
$ws.Range("L64").Value = 4656.875
$ws.Range("M64").Value = -3552
$ws.Range("N64").Value = -5152.875

$ws.Range("H67").Value = 4180.8335
$ws.Range("I67").Value = 3800
$ws.Range("J67").Value = 4656.875
$ws.Range("K67").Value = 3800
$ws.Range("L67").Value = 4656.875
$ws.Range("M67").Value = -2942
$ws.Range("N67").Value = -6372.875

$ws.Range("H76").Value = 3468.9
$ws.Range("I76").Value = 3476.5557
$ws.Range("K76").Value = 3476.5557
$ws.Range("M76").Value = -3161.5557

$ws.Range("H79").Value = 3468.9
$ws.Range("I79").Value = 3476.5557
$ws.Range("K79").Value = 3476.5557
$ws.Range("M79").Value = -2384.5557

$ws.Range("H125").Value = 864.8570999999999
$ws.Range("I125").Value = 730.8
$ws.Range("J125").Value = 1200
$ws.Range("K125").Value = 6577.2
$ws.Range("L125").Value = 10800
$ws.Range("M125").Value = -4117.2
$ws.Range("N125").Value = -15720

$ws.Range("H138").Value = 2978893.8
$ws.Range("I138").Value = 1865.3889
$ws.Range("J138").Value = 4389065
$ws.Range("K138").Value = 5596.1667
$ws.Range("L138").Value = 13167195
$ws.Range("M138").Value = -456.1666999999998
$ws.Range("N138").Value = -13177475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18388.768
$ws.Range("I32").Value = 18857.062
$ws.Range("J32").Value = 15058.667
$ws.Range("K32").Value = 18857.062
$ws.Range("L32").Value = 15058.667
$ws.Range("M32").Value = -18570.062
$ws.Range("N32").Value = -15632.667

$ws.Range("H45").Value = 2100.6667
$ws.Range("J45").Value = 1156.2858
$ws.Range("L45").Value = 1156.2858
$ws.Range("N45").Value = -1910.2858

$ws.Range("H61").Value = 47716170
$ws.Range("I61").Value = 83417890
$ws.Range("J61").Value = 113866
$ws.Range("K61").Value = 83417890
$ws.Range("L61").Value = 113866
$ws.Range("M61").Value = -83417678
$ws.Range("N61").Value = -114290

$ws.Range("H74").Value = 9694170
$ws.Range("I74").Value = 13945988
$ws.Range("J74").Value = 127581.25
$ws.Range("K74").Value = 13945988
$ws.Range("L74").Value = 127581.25
$ws.Range("M74").Value = -13945114
$ws.Range("N74").Value = -129329.25

$ws.Range("H77").Value = 9694170
$ws.Range("I77").Value = 13945988
$ws.Range("J77").Value = 127581.25
$ws.Range("K77").Value = 69729940
$ws.Range("L77").Value = 637906.25
$ws.Range("M77").Value = -69725572
$ws.Range("N77").Value = -646642.25

$ws.Range("H122").Value = 4117245
$ws.Range("I122").Value = 1828.1111
$ws.Range("J122").Value = 12348079
$ws.Range("K122").Value = 5484.3333
$ws.Range("L122").Value = 37044237
$ws.Range("M122").Value = -3034.3333
$ws.Range("N122").Value = -37049137

$ws.Range("H136").Value = 47716170
$ws.Range("I136").Value = 83417890
$ws.Range("J136").Value = 113866
$ws.Range("K136").Value = 250253670
$ws.Range("L136").Value = 341598
$ws.Range("M136").Value = -250251120
$ws.Range("N136").Value = -346698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2450.2554
$ws.Range("I31").Value = 1352.5358
$ws.Range("J31").Value = 4067.9473
$ws.Range("K31").Value = 1352.5358
$ws.Range("L31").Value = 4067.9473
$ws.Range("M31").Value = -1057.5358
$ws.Range("N31").Value = -4657.9473

$ws.Range("H34").Value = 2450.2554
$ws.Range("I34").Value = 1352.5358
$ws.Range("J34").Value = 4067.9473
$ws.Range("K34").Value = 1352.5358
$ws.Range("L34").Value = 4067.9473
$ws.Range("M34").Value = -1150.5358
$ws.Range("N34").Value = -4471.9473

$ws.Range("H107").Value = 480.6154
$ws.Range("I107").Value = 442
$ws.Range("J107").Value = 513.7143
$ws.Range("K107").Value = 442
$ws.Range("L107").Value = 513.7143
$ws.Range("M107").Value = 1478
$ws.Range("N107").Value = -4353.7143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 992.55554
$ws.Range("J131").Value = 1071.4098
$ws.Range("L131").Value = 3214.2294
$ws.Range("N131").Value = -13294.2294

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2542.0527
$ws.Range("I122").Value = 1879.5
$ws.Range("J122").Value = 3677.8572
$ws.Range("K122").Value = 5638.5
$ws.Range("L122").Value = 11033.5716
$ws.Range("M122").Value = -3188.5
$ws.Range("N122").Value = -15933.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 216.1875
$ws.Range("I55").Value = 127.5
$ws.Range("J55").Value = 304.875
$ws.Range("K55").Value = 127.5
$ws.Range("L55").Value = 304.875
$ws.Range("M55").Value = 45.5
$ws.Range("N55").Value = -650.875

$ws.Range("H68").Value = 767.64703
$ws.Range("I68").Value = 690.625
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 690.625
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = 58.375
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 767.64703
$ws.Range("I71").Value = 690.625
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 3453.125
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = 290.875
$ws.Range("N71").Value = -17488

$ws.Range("H136").Value = 159627.16
$ws.Range("I136").Value = 168492.17
$ws.Range("J136").Value = 152028.58
$ws.Range("K136").Value = 505476.51
$ws.Range("L136").Value = 456085.74
$ws.Range("M136").Value = -502926.51
$ws.Range("N136").Value = -461185.74

$ws.Range("H141").Value = 32449.5
$ws.Range("J141").Value = 32449.5
$ws.Range("L141").Value = 32449.5
$ws.Range("N141").Value = -42809.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 33483.832
$ws.Range("J80").Value = 38180.6
$ws.Range("L80").Value = 38180.6
$ws.Range("N80").Value = -40176.6

$ws.Range("H83").Value = 33483.832
$ws.Range("J83").Value = 38180.6
$ws.Range("L83").Value = 114541.8
$ws.Range("N83").Value = -124525.8

$ws.Range("H113").Value = 638.0625
$ws.Range("I113").Value = 458.14285
$ws.Range("J113").Value = 981.5454999999999
$ws.Range("K113").Value = 1374.42855
$ws.Range("L113").Value = 2944.6365
$ws.Range("M113").Value = 795.5714499999999
$ws.Range("N113").Value = -7284.6365

$ws.Range("H126").Value = 907.8333
$ws.Range("I126").Value = 907.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2723.4999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -253.4998999999998
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 58656.484
$ws.Range("I132").Value = 53694.684
$ws.Range("J132").Value = 64548.625
$ws.Range("K132").Value = 161084.052
$ws.Range("L132").Value = 193645.875
$ws.Range("M132").Value = -158554.052
$ws.Range("N132").Value = -198705.875

$ws.Range("H140").Value = 45931.668
$ws.Range("J140").Value = 45931.668
$ws.Range("L140").Value = 45931.668
$ws.Range("N140").Value = -56291.668
